$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Row 5: add a new value in column R (no pre-existing style needed)
# ---------------------------------------------------------------------------
$ws.Range("R5").Value = 5

# ---------------------------------------------------------------------------
# Row 8: fill in previously empty cells (existing styles are kept automatically)
# and create brand-new cells that need specific formatting copied in first.
# ---------------------------------------------------------------------------
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 5

$ws.Range("F8").Value = 5

$ws.Range("I19").Copy()
$ws.Range("G8").PasteSpecial($xlPasteFormats)
$ws.Range("G8").Value = 5

$ws.Range("I19").Copy()
$ws.Range("H8").PasteSpecial($xlPasteFormats)
$ws.Range("H8").Value = 5

$ws.Range("L8").Value = 4

# ---------------------------------------------------------------------------
# Row 19: fill in previously empty cells and create new cells with formatting
# ---------------------------------------------------------------------------
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 5

$ws.Range("L9").Copy()
$ws.Range("F19").PasteSpecial($xlPasteFormats)
$ws.Range("F19").Value = 5

$ws.Range("I19").Value = 5

$ws.Range("I19").Copy()
$ws.Range("J19").PasteSpecial($xlPasteFormats)
$ws.Range("J19").Value = 5

# ---------------------------------------------------------------------------
# Row 21: add a new value in column R
# ---------------------------------------------------------------------------
$ws.Range("R21").Value = 5

# ---------------------------------------------------------------------------
# Row 23: new formatted cell plus a new value in column R
# ---------------------------------------------------------------------------
$ws.Range("I19").Copy()
$ws.Range("I23").PasteSpecial($xlPasteFormats)
$ws.Range("I23").Value = 5

$ws.Range("R23").Value = 5

# ---------------------------------------------------------------------------
# Row 30: fill in previously empty cells and create several new cells
# ---------------------------------------------------------------------------
$ws.Range("C30").Value = 5
$ws.Range("D30").Value = 5
$ws.Range("E30").Value = 5

$ws.Range("I14").Copy()
$ws.Range("F30").PasteSpecial($xlPasteFormats)
$ws.Range("F30").Value = 5

$ws.Range("L5").Copy()
$ws.Range("G30").PasteSpecial($xlPasteFormats)
$ws.Range("G30").Value = 5

$ws.Range("L5").Copy()
$ws.Range("H30").PasteSpecial($xlPasteFormats)
$ws.Range("H30").Value = 5

# I30, J30 and K30 only receive the formatting (style 14) - no values
$ws.Range("L5").Copy()
$ws.Range("I30:K30").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# Row 31: fill in a previously empty cell and create two new formatted cells
# ---------------------------------------------------------------------------
$ws.Range("C31").Value = 5

$ws.Range("I19").Copy()
$ws.Range("G31").PasteSpecial($xlPasteFormats)
$ws.Range("G31").Value = 5

$ws.Range("I19").Copy()
$ws.Range("H31").PasteSpecial($xlPasteFormats)
$ws.Range("H31").Value = 5

# ---------------------------------------------------------------------------
# Clear clipboard marker and force recalculation so cached formula results
# (P4, P8, etc.) are refreshed to match the new underlying values.
# ---------------------------------------------------------------------------
$excel.CutCopyMode = $false
$excel.CalculateFull()

# ---------------------------------------------------------------------------
# Restore the final selection exactly as recorded (bottomRight pane, I30)
# ---------------------------------------------------------------------------
$ws.Range("I30").Select()
